# The "Reviews" sheet stores Date_of_scrapping as plain text (e.g. "5/2/2024").
# Update that text value to "5/7/2024" for every data row (rows 2-17),
# forcing a Text number format first so Excel does not auto-convert the
# string into a date serial value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reviews")

$dateRange = $ws.Range("A2:A17")
$dateRange.NumberFormat = "@"

for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = "5/7/2024"
}
